$d = $word.ActiveDocument

# 1. "Q-Q plot" -> "histogram" in sentence b.
$d.Content.Find.Execute(
    "We have a small sample (n = 20), so we need to create a Q-Q plot to be sure",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We have a small sample (n = 20), so we need to create a histogram to be sure",
    2)

# 2. "the distribution of sample means is normal:" -> split into two sentences
$d.Content.Find.Execute(
    "the distribution of sample means is normal:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the distribution of differences is normal in order to be able to claim that the distribution sample means is normal:",
    2)

# 3. "The data appear to be normal" -> "The differences appear to be normal"
$d.Content.Find.Execute(
    "The data appear to be normal",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The differences appear to be normal",
    2)
